$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Statusbericht 6 (19.08.) : fill in week-8 (row 8) actuals -------------
# SOLL (B8) stays driven by the existing A8 formula; IST columns get the
# new actuals and the Gesamt (sum) formulas that mirror the rows above.
$ws.Range("B8").Value = 1169.44
$ws.Range("C8").Formula = "=A8+B8"

$ws.Range("E8").Value = 1443.92
$ws.Range("F8").Value = 1169.44
$ws.Range("G8").Formula = "=E8+F8"

# --- Projektzeit: Woche 6 Fortschritt (E22) --------------------------------
$ws.Range("E22").Value = 65

# --- restore the selection Excel leaves behind after entering this data ---
[void]$ws.Range("E9").Select()
